# Fill A1:A9 with the interval sequence A001..A009 and leave the
# selection on the cell right after the last entry (B9), matching the
# "modify for interval & last" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("A001", "A002", "A003", "A004", "A005", "A006", "A007", "A008", "A009")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Range("A$row").Value = $values[$i]
}

$ws.Range("B9").Select()
